$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update Scale (E), Width (F), Height (G) for all data rows (2-45)
$ws.Range("E2:E45").Value = 8.05
$ws.Range("F2:F45").Value = 2048
$ws.Range("G2:G45").Value = 2048

# Update the view/selection to match the saved workbook state
$ws.Range("O29").Select()

$wb.Save()
